# Add "test cases for add new tasks" to the tasks sheet (sheet4) and
# switch the active tab from "deals" to "tasks".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tasks")

# ---------------------------------------------------------------------
# Fill in the task data. Values are written in the same column-major
# order the workbook's shared-string table was originally built in, so
# that newly interned strings land at the expected indices:
#   1) headers for the columns that are brand new (B, D, F, G, J, K, L)
#   2) row-2/row-3 data, column by column, left to right (A..L, N);
#      column H ("case") is itself a brand-new column, so its header is
#      written right before its data in this same pass
#   3) column M ("keyCompany") header + data, added last
# ---------------------------------------------------------------------

# Step 1 - new column headers (row 1), skipping columns that keep an
# already-existing header (C=status, E=type, I=tags, N=identifier) and
# skipping H/M which are introduced later.
$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "autoExtend"
$ws.Range("C1").Value = "status"
$ws.Range("D1").Value = "completion"
$ws.Range("E1").Value = "type"
$ws.Range("F1").Value = "priority"
$ws.Range("G1").Value = "deal"
$ws.Range("I1").Value = "tags"
$ws.Range("J1").Value = "description"
$ws.Range("K1").Value = "ownerAssignedTo"
$ws.Range("L1").Value = "keyContact"
$ws.Range("N1").Value = "identifier"

# Step 2 - row 2 / row 3 data, column by column.
$ws.Range("A2").Value = "Test Title - Task 1"
$ws.Range("A3").Value = "Test Title - Task 2"

$ws.Range("B2").Value = "Extend deadline by 14 days"
$ws.Range("B3").Value = "Extend deadline by 30 days"

$ws.Range("C2").Value = "Open"
$ws.Range("C3").Value = "Open"

$ws.Range("D2").Value = 80
$ws.Range("D3").Value = 90

$ws.Range("E2").Value = "Meeting"
$ws.Range("E3").Value = "Training"

$ws.Range("F2").Value = "High"
$ws.Range("F3").Value = "Normal"

$ws.Range("G2").Value = "Test deal -1"
$ws.Range("G3").Value = "Test deal -2"

$ws.Range("H1").Value = "case"
$ws.Range("H2").Value = "Test case-1"
$ws.Range("H3").Value = "Test case-2"

$ws.Range("I2").Value = "Test tags -1 "
$ws.Range("I3").Value = "Test tags -12"

$ws.Range("J2").Value = "Test desc -1 "
$ws.Range("J3").Value = "Test desc -2 "

$ws.Range("K2").Value = "Tejas niturkar"
$ws.Range("K3").Value = "Tejas niturkar"

$ws.Range("L2").Value = "Test contact - 1"
$ws.Range("L3").Value = "Test contact - 2"

$ws.Range("N2").Value = "test identifier - 1"
$ws.Range("N3").Value = "test identifier - 2"

# Step 3 - column M ("keyCompany") added last.
$ws.Range("M1").Value = "keyCompany"
$ws.Range("M2").Value = "Test company - 1"
$ws.Range("M3").Value = "Test company - 2"

# ---------------------------------------------------------------------
# Formatting: column D ("completion") keeps its numbers as text-like
# cells (same custom text number format used elsewhere in the workbook).
# ---------------------------------------------------------------------
$ws.Range("D1:D3").NumberFormat = "@"

# Column widths (best-fit sizing done previously by Excel for this sheet).
function Set-ExactColumnWidth($sheet, $col, $targetXmlWidth) {
    $px = [Math]::Round($targetXmlWidth * 6)
    $sheet.Columns.Item($col).ColumnWidth = ($px - 5) / 6
}

Set-ExactColumnWidth $ws 1  16.140625
Set-ExactColumnWidth $ws 2  25.28515625
Set-ExactColumnWidth $ws 4  11.140625
Set-ExactColumnWidth $ws 5  8.42578125
Set-ExactColumnWidth $ws 7  11.140625
Set-ExactColumnWidth $ws 8  10.7109375
Set-ExactColumnWidth $ws 9  11.85546875
Set-ExactColumnWidth $ws 10 11.7109375
Set-ExactColumnWidth $ws 11 17
Set-ExactColumnWidth $ws 12 14.28515625
Set-ExactColumnWidth $ws 13 15.85546875
Set-ExactColumnWidth $ws 14 16

# ---------------------------------------------------------------------
# Tab / selection changes: the "tasks" sheet becomes the active tab
# (previously "deals" was selected), with the selection left at M9.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("M9").Select()

Write-Host "tasks sheet populated"
